# Remove the row for "MAJAJD" / "Ines" (row 101) from the assignments
# table. Deleting the row shifts every subsequent row up by one and the
# table/worksheet ranges shrink accordingly (handled automatically by
# Excel when a row inside a table is deleted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("101").Delete()

# Leave the selection where the edit ended up, matching the user's
# workflow after removing the row.
$ws.Range("H108").Select()
